$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Constant" (row 5) and "r2_adj" (row 6) rows entirely -
# they no longer exist in the updated table.
$ws.Rows("5:6").Delete()

# Column header D1: "$\pi$" -> "C"
$ws.Range("D1").Value = "C"

# Row label A4: "$\pi$ Lag" -> "C Lag"
$ws.Range("A4").Value = "C Lag"

# Helper: write a numeric-looking value as genuine text (matching the
# workbook's existing convention of storing these coefficients as shared
# strings) without disturbing any cell's style/number format. We do this
# by building the text with TEXT(), copying it, and pasting values-only
# into the destination - PasteSpecial(values) carries over the text
# result but none of the formatting/style of the scratch cell.
function Set-TextValue($ref, $text) {
    $ws.Range("Z100").Formula = '="' + $text + '"'
    $ws.Range("Z100").Copy()
    $ws.Range($ref).PasteSpecial(-4163)
    $ws.Range("Z100").Clear()
}

# Row 2 (U Lag)
Set-TextValue "B2" "1.287*"
Set-TextValue "C2" "0.013"
Set-TextValue "D2" "-0.738"

# Row 3 (A Lag)
Set-TextValue "B3" "3.937"
Set-TextValue "C3" "-0.059"
Set-TextValue "D3" "44.228***"

# Row 4 (C Lag)
Set-TextValue "B4" "-0.343"
Set-TextValue "C4" "-0.001"
Set-TextValue "D4" "-0.615***"
